# Fruta / hortaliza, semanal
# Inserts 4 new weekly price rows (new report date 2021-10-05, serial 44474)
# before the existing row 117, shifting the remaining historical rows down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 117:134 down to 121:138, inserting 4 blank rows at 117.
$ws.Rows("117:120").Insert()

# --- New row 117 ---
$ws.Range("A117").Value = 6
$ws.Range("B117").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C117").Value = "Metropolitana"
$ws.Range("D117").Value = 44474
$ws.Range("E117").Value = 13
$ws.Range("F117").Value = 100112026
$ws.Range("G117").Value = "Haba"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 800
$ws.Range("K117").Value = 6000
$ws.Range("L117").Value = 7000
$ws.Range("M117").Value = 6562
$ws.Range("N117").Value = "$/saco 25 kilos"
$ws.Range("O117").Value = "Región Metropolitana"
$ws.Range("P117").Value = 262
$ws.Range("Q117").Value = 25
$ws.Range("R117").Value = "Hortaliza"

# --- New row 118 ---
$ws.Range("A118").Value = 6
$ws.Range("B118").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C118").Value = "Metropolitana"
$ws.Range("D118").Value = 44474
$ws.Range("E118").Value = 13
$ws.Range("F118").Value = 100112026
$ws.Range("G118").Value = "Haba"
$ws.Range("H118").Value = "Sin especificar"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 200
$ws.Range("K118").Value = 6000
$ws.Range("L118").Value = 6000
$ws.Range("M118").Value = 6000
$ws.Range("N118").Value = "$/saco 25 kilos"
$ws.Range("O118").Value = "Región de Coquimbo"
$ws.Range("P118").Value = 240
$ws.Range("Q118").Value = 25
$ws.Range("R118").Value = "Hortaliza"

# --- New row 119 ---
$ws.Range("A119").Value = 6
$ws.Range("B119").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C119").Value = "Metropolitana"
$ws.Range("D119").Value = 44474
$ws.Range("E119").Value = 13
$ws.Range("F119").Value = 100112026
$ws.Range("G119").Value = "Haba"
$ws.Range("H119").Value = "Sin especificar"
$ws.Range("I119").Value = "Segunda"
$ws.Range("J119").Value = 250
$ws.Range("K119").Value = 5000
$ws.Range("L119").Value = 5000
$ws.Range("M119").Value = 5000
$ws.Range("N119").Value = "$/saco 25 kilos"
$ws.Range("O119").Value = "Región Metropolitana"
$ws.Range("P119").Value = 200
$ws.Range("Q119").Value = 25
$ws.Range("R119").Value = "Hortaliza"

# --- New row 120 ---
$ws.Range("A120").Value = 6
$ws.Range("B120").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C120").Value = "Metropolitana"
$ws.Range("D120").Value = 44474
$ws.Range("E120").Value = 13
$ws.Range("F120").Value = 100112026
$ws.Range("G120").Value = "Haba"
$ws.Range("H120").Value = "Sin especificar"
$ws.Range("I120").Value = "Segunda"
$ws.Range("J120").Value = 100
$ws.Range("K120").Value = 5000
$ws.Range("L120").Value = 5000
$ws.Range("M120").Value = 5000
$ws.Range("N120").Value = "$/saco 25 kilos"
$ws.Range("O120").Value = "Región de Coquimbo"
$ws.Range("P120").Value = 200
$ws.Range("Q120").Value = 25
$ws.Range("R120").Value = "Hortaliza"
